$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns for the crypto list.
# D-column numeric-looking values must stay as text, matching the source
# data (inline/shared strings), so NumberFormat is forced to Text ("@")
# before the write and the cell style is reset to "Normal" afterward so
# no new style index is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.887.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.913.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.92%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.913.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("E11").Value = "  -5.34%  "
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.92%  "
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.395.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.835.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.914.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0875"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.69%  "
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.91%  "
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "379.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.686.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -2.57%  "
